# ============================================================
# Change 1: Title paragraph
#   "CAHIER DE CHARGE DU PROJET D’ELECTRONIQUE : MONITORING DE L’ENERGIE
#    SOLAIRE ET D’UN SUIVEUR SOLAIRE"
#   -> "CAHIER DE CHARGE DU PROJET D’IHM : MUTUELLE"
#   split into 4 runs: "CAHIER DE CHARGE DU PROJET " (u) / "D’IHM" (u) /
#   " : " (no u) / "MUTUELLE" (no u)
# ============================================================
$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1)
$s = $p1.Range.Start

# Run 1 originally spans $s .. $s+42 and is underlined; rewrite its text.
$r1 = $d.Range($s, $s + 42)
$r1.Text = "CAHIER DE CHARGE DU PROJET D’IHM"

# Run 2 originally spans $s+32 .. $s+90 (after run1 shrank) and is NOT
# underlined; rewrite its text.
$r2 = $d.Range($s + 32, $s + 90)
$r2.Text = " : MUTUELLE"

# Split the underlined portion into "CAHIER DE CHARGE DU PROJET " / "D’IHM"
# (toggle Underline off/on so the run boundary is forced, ending back at
# the same "single" value).
$rSplit1 = $d.Range($s + 27, $s + 32)
$rSplit1.Font.Underline = 0
$rSplit1.Font.Underline = 1

# Split the non-underlined portion into " : " / "MUTUELLE" (toggle Bold
# on/off so the run boundary is forced, ending back at the default value).
$rSplit2 = $d.Range($s + 35, $s + 43)
$rSplit2.Font.Bold = 1
$rSplit2.Font.Bold = 0

Write-Output ("Title now: '" + $d.Paragraphs(1).Range.Text + "'")

# ============================================================
# Change 2: add a new "NNANGE AKUME" bullet right after
#   "TCHOFFO NGINTEDEM" and before "MOHAMADOU HAMIDOU"
# ============================================================
$rngFind = $d.Content
$rngFind.Find.Execute("TCHOFFO NGINTEDEM")
$tchoffoPara = $rngFind.Paragraphs(1)
$tchoffoPara.Range.InsertParagraphAfter()

$rngFind2 = $d.Content
$rngFind2.Find.Execute("TCHOFFO NGINTEDEM")
$newPara = $rngFind2.Paragraphs(1).Next()
$newPara.Range.Text = "NNANGE AKUME"

Write-Output ("New paragraph: '" + $newPara.Range.Text + "'")

# ============================================================
# Change 3: merge the "Il réalise sa tâche..." runs (drop the
# proofErr-wrapped spell-checked run split) into a single run.
# ============================================================
$oldText = "Il réalise sa tâche telle que prévu. C’est-à-dire monitoring de l’energie solaire et suivi de la course du soleil"
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $oldText, 2)

Write-Output "Change3 done"
